{"js": "// Replace the date header and each division problem's text with the\n// updated values described by the diff. Every search string below is\n// unique in the document, so a simple search + insertText(\"Replace\")\n// per pair is sufficient and safe.\nconst replacements = [\n  [\"2025-07-16 Wednesday\", \"2025-07-17 Thursday\"],\n  [\"44\u00f75=\", \"22\u00f76=\"],\n  [\"28\u00f77=\", \"84\u00f73=\"],\n  [\"36\u00f72=\", \"24\u00f72=\"],\n  [\"40\u00f78=\", \"42\u00f78=\"],\n  [\"70\u00f73=\", \"42\u00f76=\"],\n  [\"89\u00f77=\", \"58\u00f72=\"],\n  [\"29\u00f74=\", \"45\u00f75=\"],\n  [\"43\u00f76=\", \"47\u00f75=\"],\n  [\"84\u00f78=\", \"18\u00f78=\"],\n  [\"56\u00f79=\", \"23\u00f79=\"],\n  [\"49\u00f72=\", \"60\u00f78=\"],\n  [\"49\u00f78=\", \"42\u00f73=\"],\n  [\"68\u00f79=\", \"33\u00f72=\"],\n  [\"89\u00f72=\", \"87\u00f74=\"],\n  [\"54\u00f72=\", \"72\u00f72=\"],\n  [\"63\u00f73=\", \"75\u00f75=\"],\n  [\"42\u00f74=\", \"24\u00f72=\"],\n  [\"50\u00f77=\", \"22\u00f75=\"],\n  [\"27\u00f72=\", \"13\u00f74=\"],\n  [\"31\u00f74=\", \"98\u00f74=\"],\n  [\"20\u00f78=\", \"25\u00f75=\"],\n  [\"13\u00f77=\", \"73\u00f79=\"],\n  [\"70\u00f79=\", \"94\u00f74=\"],\n  [\"52\u00f78=\", \"51\u00f77=\"],\n  [\"67\u00f75=\", \"84\u00f72=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date header and each division problem's text with the\n# updated values described by the diff. Every \"find\" string below is\n# unique in the document, so a plain Find.Execute replace-all per pair\n# is safe and deterministic.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-07-16 Wednesday\", \"2025-07-17 Thursday\"),\n    @(\"44\u00f75=\", \"22\u00f76=\"),\n    @(\"28\u00f77=\", \"84\u00f73=\"),\n    @(\"36\u00f72=\", \"24\u00f72=\"),\n    @(\"40\u00f78=\", \"42\u00f78=\"),\n    @(\"70\u00f73=\", \"42\u00f76=\"),\n    @(\"89\u00f77=\", \"58\u00f72=\"),\n    @(\"29\u00f74=\", \"45\u00f75=\"),\n    @(\"43\u00f76=\", \"47\u00f75=\"),\n    @(\"84\u00f78=\", \"18\u00f78=\"),\n    @(\"56\u00f79=\", \"23\u00f79=\"),\n    @(\"49\u00f72=\", \"60\u00f78=\"),\n    @(\"49\u00f78=\", \"42\u00f73=\"),\n    @(\"68\u00f79=\", \"33\u00f72=\"),\n    @(\"89\u00f72=\", \"87\u00f74=\"),\n    @(\"54\u00f72=\", \"72\u00f72=\"),\n    @(\"63\u00f73=\", \"75\u00f75=\"),\n    @(\"42\u00f74=\", \"24\u00f72=\"),\n    @(\"50\u00f77=\", \"22\u00f75=\"),\n    @(\"27\u00f72=\", \"13\u00f74=\"),\n    @(\"31\u00f74=\", \"98\u00f74=\"),\n    @(\"20\u00f78=\", \"25\u00f75=\"),\n    @(\"13\u00f77=\", \"73\u00f79=\"),\n    @(\"70\u00f79=\", \"94\u00f74=\"),\n    @(\"52\u00f78=\", \"51\u00f77=\"),\n    @(\"67\u00f75=\", \"84\u00f72=\")\n)\n\nforeach ($pair in $pairs) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n"}
